# Busqueda en Wikipedia por datos Excel
# Re-runs the "search" over the Descripcion column (A) and writes the
# first Wikipedia result into PrimerResultado (B). For this data set every
# search term resolved to itself (exact-title matches), so column B mirrors
# column A for each processed row.
#
# Row 5 ("Car") is looked up/written before row 4 ("Driver") -- the search
# results come back out of request order -- which is why "Car" claims the
# lower shared-string slot even though it ends up a row below "Driver" in
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 2-3: overwrite the old "hola" lookups with the new search
# terms / results.
$ws.Range("A2").Value = "Montana"
$ws.Range("B2").Value = "Montana"

$ws.Range("A3").Value = "Peru"
$ws.Range("B3").Value = "Peru"

# New row 5 ("Car") is resolved first...
$ws.Range("A5").Value = "Car"
$ws.Range("B5").Value = "Car"
$ws.Range("A5").Locked = $false

# ...then new row 4 ("Driver") is appended.
$ws.Range("A4").Value = "Driver"
$ws.Range("B4").Value = "Driver"
$ws.Range("C4").Locked = $false

# Page is set up for printing the refreshed results.
$ws.PageSetup.Orientation = 1

# Leave the cursor where the next search term would be typed.
$ws.Range("B4").Select() | Out-Null
